$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.94"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.01"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.363"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05853"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.394"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.372"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8135"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.018"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1419"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04312"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCXBestin24h"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07334"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.205"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09391"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001598"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04825"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005890"
$ws.Range("E18").Value = "17OneONEWorstin24h"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006007"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004084"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009814"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.710"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002471"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003000"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.006329"
$ws.Range("E43").Value = "42KickTokenKICK"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005074"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005628"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7610"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.09025"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
